$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$ws1 = $wb.Worksheets.Item("Rushing")

# A.Jones
$ws1.Range("C2").Value = 5
$ws1.Range("D2").Value = 6

# A.Dillon
$ws1.Range("C3").Value = 100
$ws1.Range("D3").Value = 51

# K.Hill
$ws1.Range("C4").Value = 76
$ws1.Range("D4").Value = 51

# A.Lazard
$ws1.Range("C7").Value = 1

# --- Receiving sheet ---
$ws2 = $wb.Worksheets.Item("Receiving")

# A.Dillon
$ws2.Range("C2").Value = 52
$ws2.Range("D2").Value = 43

# K.Hill
$ws2.Range("C3").Value = 23
$ws2.Range("D3").Value = 19

# D.Adams
$ws2.Range("C5").Value = 109
$ws2.Range("D5").Value = 91
$ws2.Range("E5").Value = 37

# M.Valdes-Scantling
$ws2.Range("C6").Value = 19
$ws2.Range("D6").Value = 15
$ws2.Range("E6").Value = 22
$ws2.Range("F6").Value = 6

# A.Lazard
$ws2.Range("C7").Value = 41
$ws2.Range("D7").Value = 33

# J.Winfree
$ws2.Range("C11").Value = 7

# T.Davis
$ws2.Range("E12").Value = 1
$ws2.Range("F12").Value = 1

# M.Lewis
$ws2.Range("C13").Value = 23
$ws2.Range("D13").Value = 21
$ws2.Range("E13").Value = 4
$ws2.Range("F13").Value = 3

# J.Deguara
$ws2.Range("C14").Value = 19
$ws2.Range("D14").Value = 15
